$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D="43.173.90"; E="  +0.26%  "},
    @{Row=3; D="2.321.90"; E="  +0.83%  "},
    @{Row=4; D=$null; E="  -0.01%  "},
    @{Row=5; D=$null; E="  +0.12%  "},
    @{Row=6; D="99.40"; E="  -0.14%  "},
    @{Row=7; D="0.507"; E="  +0.18%  "},
    @{Row=8; D=$null; E="  +0.01%  "},
    @{Row=9; D=$null; E="  +1.66%  "},
    @{Row=10; D="36.28"; E="  +5.41%  "},
    @{Row=11; D=$null; E="  -0.71%  "},
    @{Row=13; D="17.57"; E="  -0.78%  "},
    @{Row=14; D="6.92"; E="  +1.75%  "},
    @{Row=15; D="2.684.05"; E="  +0.84%  "},
    @{Row=16; D="2.375.45"; E="  +2.11%  "},
    @{Row=17; D="0.797"; E="  -1.33%  "},
    @{Row=18; D="43.101.29"; E="  +0.34%  "},
    @{Row=19; D="12.85"; E="  +4.27%  "},
    @{Row=20; D=$null; E="  +1.83%  "},
    @{Row=21; D=$null; E="  +0.37%  "},
    @{Row=22; D="68.19"; E="  +0.49%  "},
    @{Row=23; D="240.52"; E="  +1.47%  "},
    @{Row=24; D=$null; E="  -1.74%  "},
    @{Row=25; D="2.45"; E="  -0.85%  "},
    @{Row=26; D=$null; E="  -0.11%  "},
    @{Row=27; D="25.48"; E="  +3.78%  "},
    @{Row=28; D="168.50"; E="  -0.17%  "},
    @{Row=29; D="34.28"; E="  +1.16%  "},
    @{Row=30; D="9.20"; E="  +0.16%  "},
    @{Row=31; D=$null; E="  -2.48%  "},
    @{Row=32; D=$null; E="  +3.22%  "},
    @{Row=33; D="1.00"; E="  +0.00%  "},
    @{Row=34; D="4.74"; E="  +3.86%  "},
    @{Row=35; D=$null; E="  +3.95%  "},
    @{Row=36; D=$null; E="  -1.14%  "},
    @{Row=37; D="0.0695"; E="  -0.64%  "},
    @{Row=38; D=$null; E="  +0.29%  "},
    @{Row=39; D=$null; E="  +0.40%  "},
    @{Row=40; D=$null; E="  -1.95%  "},
    @{Row=41; D=$null; E="  +0.28%  "},
    @{Row=42; D="2.002.80"; E="  +0.08%  "},
    @{Row=43; D="0.0290"; E="  +1.37%  "},
    @{Row=44; D=$null; E="  -4.27%  "},
    @{Row=45; D="10.12"; E="  -0.10%  "},
    @{Row=46; D="17.61"; E="  -0.81%  "},
    @{Row=47; D="2.88"; E="  +0.35%  "},
    @{Row=48; D="54.94"; E="  -1.34%  "},
    @{Row=49; D="75.72"; E="  +7.81%  "},
    @{Row=50; D="2.549.23"; E="  +0.82%  "},
    @{Row=51; D="1.54"; E="  +1.88%  "}
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.D -ne $null) {
        $dcell = $ws.Cells.Item($r, 4)
        $dcell.NumberFormat = "@"
        $dcell.Value = $item.D
        $dcell.Style = "Normal"
    }
    $ecell = $ws.Cells.Item($r, 5)
    $ecell.Value = $item.E
}
